# Resume update: "Graduating July 2019" -> "Graduated July 2019"
#
# The canonical diff shows the single run containing "Graduating July 2019"
# split into three runs -- "Graduat", "ed", " July 2019" -- with the
# document's "_GoBack" bookmark (which Word drops at the site of the most
# recent edit) relocated from its old home (end of the "Match parking
# rates..." bullet) to sit right after the newly retyped "ed".
#
# We reproduce that run layout explicitly: temporarily bracket the "ing"
# text with two bookmarks (forcing Word to split the run at those two
# points), overwrite just that bracketed "ing" text with "ed", then drop
# the temporary left-hand bookmark and keep/rename the right-hand one as
# "_GoBack" (adding a bookmark named "_GoBack" automatically relocates the
# document's unique _GoBack bookmark, removing it from its previous spot).

$d = $word.ActiveDocument

$r = $d.Content
$found = $r.Find.Execute("Graduating July 2019", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Graduating July 2019'"
}

$phraseStart = $r.Start

# Offsets within "Graduating July 2019":
#   "Graduat" = 7 chars, "ing" = 3 chars, then " July 2019" follows.
$leftSplit = $phraseStart + 7
$rightSplit = $phraseStart + 10

# Bracket "ing" with two bookmarks so the run gets split on both sides of it.
$d.Bookmarks.Add("ZZTempGradSplit", $d.Range($leftSplit, $leftSplit))
$d.Bookmarks.Add("_GoBack", $d.Range($rightSplit, $rightSplit))

# Overwrite just the bracketed "ing" text with "ed" (stays within the two
# bookmark-induced run boundaries, so it lands in its own run).
$ingRange = $d.Range($leftSplit, $rightSplit)
$ingRange.Text = "ed"

# Drop the temporary left marker; the real "_GoBack" bookmark (and the run
# split it caused) remains right after "ed".
$d.Bookmarks("ZZTempGradSplit").Delete()
